$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 305.04166
$ws.Range("J17").Value = 314.9091
$ws.Range("L17").Value = 944.7273
$ws.Range("N17").Value = -1280.7273
$ws.Range("H86").Value = 1188.6
$ws.Range("I86").Value = 1265.1111
$ws.Range("K86").Value = 1265.1111
$ws.Range("M86").Value = -142.1111000000001
$ws.Range("H89").Value = 1188.6
$ws.Range("I89").Value = 1265.1111
$ws.Range("K89").Value = 6325.5555
$ws.Range("M89").Value = -709.5555000000004
$ws.Range("H92").Value = 656.7917
$ws.Range("I92").Value = 520.5625
$ws.Range("K92").Value = 520.5625
$ws.Range("M92").Value = 727.4375
$ws.Range("H105").Value = 34900
$ws.Range("J105").Value = 34900
$ws.Range("L105").Value = 34900
$ws.Range("N105").Value = -41888
$ws.Range("H107").Value = 487.5862
$ws.Range("I107").Value = 360.08694
$ws.Range("K107").Value = 360.08694
$ws.Range("M107").Value = 1559.91306
$ws.Range("H111").Value = 1246.75
$ws.Range("I111").Value = 662.3333
$ws.Range("K111").Value = 1986.9999
$ws.Range("M111").Value = 1080.0001
$ws.Range("H137").Value = 2860.6843
$ws.Range("I137").Value = 2333.7334
$ws.Range("J137").Value = 4836.75
$ws.Range("K137").Value = 7001.2002
$ws.Range("L137").Value = 14510.25
$ws.Range("M137").Value = -4451.2002
$ws.Range("N137").Value = -19610.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1827.5
$ws.Range("I2").Value = 1227.5652
$ws.Range("J2").Value = 3798.7144
$ws.Range("K2").Value = 1227.5652
$ws.Range("L2").Value = 3798.7144
$ws.Range("M2").Value = -1114.5652
$ws.Range("N2").Value = -4024.7144
$ws.Range("H76").Value = 54513.375
$ws.Range("J76").Value = 54513.375
$ws.Range("L76").Value = 54513.375
$ws.Range("N76").Value = -55189.375
$ws.Range("H79").Value = 54513.375
$ws.Range("J79").Value = 54513.375
$ws.Range("L79").Value = 54513.375
$ws.Range("N79").Value = -56853.375
$ws.Range("H116").Value = 1827.5
$ws.Range("I116").Value = 1227.5652
$ws.Range("J116").Value = 3798.7144
$ws.Range("K116").Value = 1227.5652
$ws.Range("L116").Value = 3798.7144
$ws.Range("M116").Value = 1066.4348
$ws.Range("N116").Value = -8386.714400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1827.5
$ws.Range("I3").Value = 1227.5652
$ws.Range("J3").Value = 3798.7144
$ws.Range("K3").Value = 1227.5652
$ws.Range("L3").Value = 3798.7144
$ws.Range("M3").Value = -1113.5652
$ws.Range("N3").Value = -4026.7144
$ws.Range("H134").Value = 1665.5555
$ws.Range("I134").Value = 1334.0952
$ws.Range("K134").Value = 4002.2856
$ws.Range("M134").Value = -1467.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2365.9656
$ws.Range("I16").Value = 2716.1052
$ws.Range("K16").Value = 2716.1052
$ws.Range("M16").Value = -2429.1052
$ws.Range("H31").Value = 3313.2778
$ws.Range("I31").Value = 1974.2142
$ws.Range("K31").Value = 1974.2142
$ws.Range("M31").Value = -1679.2142
$ws.Range("H34").Value = 3313.2778
$ws.Range("I34").Value = 1974.2142
$ws.Range("K34").Value = 1974.2142
$ws.Range("M34").Value = -1772.2142
$ws.Range("H86").Value = 100006060
$ws.Range("I86").Value = 100006060
$ws.Range("K86").Value = 100006060
$ws.Range("M86").Value = -100004937
$ws.Range("H89").Value = 100006060
$ws.Range("I89").Value = 100006060
$ws.Range("K89").Value = 500030300
$ws.Range("M89").Value = -500024684
$ws.Range("H113").Value = 2365.9656
$ws.Range("I113").Value = 2716.1052
$ws.Range("K113").Value = 2716.1052
$ws.Range("M113").Value = -546.1052
$ws.Range("H132").Value = 3066.2
$ws.Range("I132").Value = 2962.6
$ws.Range("K132").Value = 8887.799999999999
$ws.Range("M132").Value = -6357.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2119.7144
$ws.Range("I3").Value = 2119.7144
$ws.Range("K3").Value = 6359.1432
$ws.Range("M3").Value = -6247.1432
$ws.Range("H5").Value = 490.06668
$ws.Range("J5").Value = 519.1667
$ws.Range("L5").Value = 1557.5001
$ws.Range("N5").Value = -1781.5001
$ws.Range("H38").Value = 1286.4706
$ws.Range("I38").Value = 2601.5
$ws.Range("K38").Value = 7804.5
$ws.Range("M38").Value = -7457.5
$ws.Range("H41").Value = 2249
$ws.Range("J41").Value = 2249
$ws.Range("L41").Value = 6747
$ws.Range("N41").Value = -7423
$ws.Range("H62").Value = 3898
$ws.Range("J62").Value = 3898
$ws.Range("L62").Value = 11694
$ws.Range("N62").Value = -13066
$ws.Range("H64").Value = 13624.75
$ws.Range("J64").Value = 12333.333
$ws.Range("L64").Value = 36999.999
$ws.Range("N64").Value = -37539.999
$ws.Range("H65").Value = 3898
$ws.Range("J65").Value = 3898
$ws.Range("L65").Value = 35082
$ws.Range("N65").Value = -41946
$ws.Range("H67").Value = 13624.75
$ws.Range("J67").Value = 12333.333
$ws.Range("L67").Value = 36999.999
$ws.Range("N67").Value = -38871.999
$ws.Range("H68").Value = 500
$ws.Range("J68").Value = 500
$ws.Range("L68").Value = 1500
$ws.Range("N68").Value = -3122
$ws.Range("H71").Value = 500
$ws.Range("J71").Value = 500
$ws.Range("L71").Value = 4500
$ws.Range("N71").Value = -12612
$ws.Range("H108").Value = 642.5
$ws.Range("I108").Value = 642.5
$ws.Range("K108").Value = 1927.5
$ws.Range("M108").Value = 952.5
$ws.Range("H114").Value = 54665.25
$ws.Range("I114").Value = 599
$ws.Range("J114").Value = 72687.336
$ws.Range("K114").Value = 1797
$ws.Range("L114").Value = 218062.008
$ws.Range("M114").Value = 1457
$ws.Range("N114").Value = -224570.008
$ws.Range("H135").Value = 490.06668
$ws.Range("J135").Value = 519.1667
$ws.Range("L135").Value = 4672.5003
$ws.Range("N135").Value = -9742.5003
$ws.Range("H138").Value = 5875.143
$ws.Range("I138").Value = 4354.3335
$ws.Range("K138").Value = 13063.0005
$ws.Range("M138").Value = -7923.000499999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4235.5386
$ws.Range("I132").Value = 4029.8572
$ws.Range("K132").Value = 12089.5716
$ws.Range("M132").Value = -9559.571599999999
$ws.Range("H134").Value = 13000
$ws.Range("J134").Value = 13000
$ws.Range("L134").Value = 39000
$ws.Range("N134").Value = -44070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 56910.35
$ws.Range("I22").Value = 129824.125
$ws.Range("J22").Value = 8301.166999999999
$ws.Range("K22").Value = 129824.125
$ws.Range("L22").Value = 8301.166999999999
$ws.Range("M22").Value = -129529.125
$ws.Range("N22").Value = -8891.166999999999
$ws.Range("H27").Value = 56910.35
$ws.Range("I27").Value = 129824.125
$ws.Range("J27").Value = 8301.166999999999
$ws.Range("K27").Value = 129824.125
$ws.Range("L27").Value = 8301.166999999999
$ws.Range("M27").Value = -129717.125
$ws.Range("N27").Value = -8515.166999999999
$ws.Range("H40").Value = 5592.25
$ws.Range("I40").Value = 5846.091
$ws.Range("J40").Value = 2800
$ws.Range("K40").Value = 5846.091
$ws.Range("L40").Value = 2800
$ws.Range("M40").Value = -5710.091
$ws.Range("N40").Value = -3072
$ws.Range("H55").Value = 1315.1538
$ws.Range("I55").Value = 389.1111
$ws.Range("K55").Value = 389.1111
$ws.Range("M55").Value = -216.1111
$ws.Range("H132").Value = 3074.9744
$ws.Range("I132").Value = 3761.05
$ws.Range("J132").Value = 2352.7896
$ws.Range("K132").Value = 11283.15
$ws.Range("L132").Value = 7058.3688
$ws.Range("M132").Value = -8753.150000000001
$ws.Range("N132").Value = -12118.3688

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 306.17392
$ws.Range("I107").Value = 308.25
$ws.Range("K107").Value = 924.75
$ws.Range("M107").Value = 995.25
$ws.Range("H113").Value = 720.6957
$ws.Range("I113").Value = 1108.1538
$ws.Range("J113").Value = 217
$ws.Range("K113").Value = 3324.4614
$ws.Range("L113").Value = 651
$ws.Range("M113").Value = -1154.4614
$ws.Range("N113").Value = -4991
$ws.Range("H122").Value = 16132504
$ws.Range("J122").Value = 2758.182
$ws.Range("L122").Value = 8274.545999999998
$ws.Range("N122").Value = -13174.546
$ws.Range("H132").Value = 1496.6333
$ws.Range("I132").Value = 1256.16
$ws.Range("K132").Value = 3768.48
$ws.Range("M132").Value = -1238.48
